# Applies the "Fixed images and instructions for README-Alexa" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new Notes column
$ws.Range("F1").Value = "Notes"

# Row 2: Person -> Adam, Status -> Done, Notes -> Split out Athena from Main README
$ws.Range("D2").Value = "Adam"
$ws.Range("E2").Value = "Done"
$ws.Range("F2").Value = "Split out Athena from Main README"

# Row 5: Status -> 50% complete, Notes -> Done form Alexa
$ws.Range("E5").Value = "50% complete"
$ws.Range("F5").Value = "Done form Alexa"

# Row 9: Status -> 50% complete, Notes -> Need to validate
$ws.Range("E9").Value = "50% complete"
$ws.Range("F9").Value = "Need to validate"

# Row 14: Status -> 50% complete, Notes -> Need to validate
$ws.Range("E14").Value = "50% complete"
$ws.Range("F14").Value = "Need to validate"

# Row 15: Status -> Done
$ws.Range("E15").Value = "Done"

# Row 16: Notes -> URL, formatted in a Menlo 9pt monospace font
$ws.Range("F16").Value = "https://github.com/voicehacks/setup-local-recommendations/blob/master/speech-assets/InteractionModel.json"
$ws.Range("F16").Font.Name = "Menlo"
$ws.Range("F16").Font.Size = 9
$ws.Range("F16").Font.Color = 0

# Row 19: Status -> 20% completed, with percent number format applied
$ws.Range("E19").Value = "20% completed"
$ws.Range("E19").NumberFormat = "0%"

# Row 26: Status -> Done
$ws.Range("E26").Value = "Done"

# Row 36: Person -> Adam
$ws.Range("D36").Value = "Adam"

# Update the view to match the committed selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D37").Select()
